# internal meeting minutes June 26 2018.docx
#
# Three changes:
#  1. Attendees cell (table 1, row 6, col 2): drop the trailing ", Kevin"
#     run and move the "_GoBack" bookmark to sit right after "Yigang"
#     (before the spellEnd proof-error marker).
#  2. "Meeting objective" table (table 3, row 2, col 2): merge the
#     "Evaluate the " and "UI development" runs into a single run,
#     leaving the following " " and "(Contact us, home page)" runs
#     untouched.
#  3. Remove the stray "_GoBack" bookmark that used to sit in the empty
#     paragraph right before the "To do list" heading.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Attendees cell -------------------------------------------------
$attendeesTable = $d.Tables(1)
$attendeesCell = $attendeesTable.Cell(6, 2)
$cellRange = $attendeesCell.Range
$attendeesRange = $d.Range($cellRange.Start, $cellRange.End)

$attendeesXml = '<w:p ' + $wNs + '>' +
  '<w:r><w:t xml:space="preserve">Wang Zhuowei, Luo Hao Nan, Yang </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Mingqi</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Gui</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Yuqi</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve">, Li </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>Yigang</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'

$attendeesRange.InsertXML($attendeesXml)

# --- 2. "Evaluate the UI development" cell ------------------------------
$objectiveTable = $d.Tables(3)
$objectiveCell = $objectiveTable.Cell(2, 2)
$objCellRange = $objectiveCell.Range
$objectiveRange = $d.Range($objCellRange.Start, $objCellRange.End)

$objectiveXml = '<w:p ' + $wNs + '>' +
  '<w:r><w:t>Evaluate the UI development</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t>(Contact us, home page)</w:t></w:r>' +
  '</w:p>'

$objectiveRange.InsertXML($objectiveXml)

# --- 3. Remove the stray "_GoBack" bookmark paragraph -------------------
$tableBefore = $d.Tables(3)
$tableAfter = $d.Tables(4)
$betweenRange = $d.Range($tableBefore.Range.End, $tableAfter.Range.Start)
$bookmarkParaRange = $betweenRange.Paragraphs(1).Range
$bookmarkRange = $d.Range($bookmarkParaRange.Start, $bookmarkParaRange.End)

$emptyParaXml = '<w:p ' + $wNs + '>' +
  '<w:pPr><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr>' +
  '</w:p>'

$bookmarkRange.InsertXML($emptyParaXml)
